$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "PhoneNumber"
$ws.Range("E1").Value = "JoinDate"

# Data row (seeded admin user)
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "admin@admin.com"
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = "2022-02-11 05:00:15.8400000"

# Hyperlink the email address (adds the Hyperlink cell style too)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:admin@admin.com")

# Column widths to fit the seeded content
$ws.Columns.Item(1).ColumnWidth = 13.67
$ws.Columns.Item(2).ColumnWidth = 11.67
$ws.Columns.Item(3).ColumnWidth = 20
$ws.Columns.Item(4).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 25.67

# Active selection left where the user last clicked
$null = $ws.Range("E8").Select()
